$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated job match data for rows 2-6
$data = @(
    @("AI/ML Engineer (Mid-Level)", "OnPoint Consulting, Inc", "Bethesda, MD, US USA", 13.3, "Data Scientist, Copilot, TensorFlow, PyTorch, Azure ML, Synapse, Data Lake, CI/CD, Git, Python", "2026-02-20", "https://www.indeed.com/viewjob?jk=ecd65ee8e23a06be"),
    @("AI Applications Development Intern", "Zoox", "Foster City, CA, US USA", 11.1, "LangChain, RAG, LLaMA, Copilot, Prompt Engineering, Git, Python, SQL, R, Java", "2026-02-21", "https://www.indeed.com/viewjob?jk=2d2402acc3891444"),
    @("Midlevel Software Engineer", "Liveworld", "Remote, US USA", 11.1, "RAG, TensorFlow, PyTorch, Git, PostgreSQL, Python, SQL, R, Java, Scala", "2026-02-20", "https://www.indeed.com/viewjob?jk=7dffdc68dbca47db"),
    @("Data Engineer", "Driven Brands", "Charlotte, NC, US USA", 10, "BigQuery, CI/CD, Snowflake, BigQuery, MySQL, SQL, R, Scala, Optimization", "2026-02-21", "https://www.indeed.com/viewjob?jk=039063e25cce7591"),
    @("Data Scientist - Quantitative Trading", "TotalEnergies", "Houston, TX, US USA", 10, "Data Scientist, Docker, CI/CD, MongoDB, Python, SQL, R, Scala, Optimization", "2026-02-21", "https://www.indeed.com/viewjob?jk=51838a211f043506")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]

    # Force the "Posted At" cell to stay plain text so the date-like
    # string isn't auto-converted into a date serial number.
    $dateCell = $ws.Cells.Item($row, 6)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item[5]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 7).Value = $item[6]
    $row++
}

# Remove the now-obsolete 7th data row entirely
$ws.Rows.Item(7).Delete()
